$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 84
$ws.Range("I33").Value = 93.5
$ws.Range("K33").Value = 93.5
$ws.Range("M33").Value = 135.5
$ws.Range("H132").Value = 2944.5264
$ws.Range("I132").Value = 3182.9412
$ws.Range("J132").Value = 918
$ws.Range("K132").Value = 9548.8236
$ws.Range("L132").Value = 2754
$ws.Range("M132").Value = -7018.8236
$ws.Range("N132").Value = -7814
$ws.Range("H137").Value = 27249.486
$ws.Range("I137").Value = 1558.5217
$ws.Range("J137").Value = 64180.25
$ws.Range("K137").Value = 4675.5651
$ws.Range("L137").Value = 192540.75
$ws.Range("M137").Value = -2125.5651
$ws.Range("N137").Value = -197640.75
$ws.Range("H138").Value = 2393.5818
$ws.Range("I138").Value = 1508.4166
$ws.Range("J138").Value = 4070.7368
$ws.Range("K138").Value = 4525.2498
$ws.Range("L138").Value = 12212.2104
$ws.Range("M138").Value = 614.7502000000004
$ws.Range("N138").Value = -22492.2104

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1370.3704
$ws.Range("I2").Value = 1207.1875
$ws.Range("J2").Value = 1607.7273
$ws.Range("K2").Value = 1207.1875
$ws.Range("L2").Value = 1607.7273
$ws.Range("M2").Value = -1094.1875
$ws.Range("N2").Value = -1833.7273
$ws.Range("H74").Value = 5383.0415
$ws.Range("I74").Value = 6128.95
$ws.Range("J74").Value = 1653.5
$ws.Range("K74").Value = 6128.95
$ws.Range("L74").Value = 1653.5
$ws.Range("M74").Value = -5254.95
$ws.Range("N74").Value = -3401.5
$ws.Range("H77").Value = 5383.0415
$ws.Range("I77").Value = 6128.95
$ws.Range("J77").Value = 1653.5
$ws.Range("K77").Value = 30644.75
$ws.Range("L77").Value = 8267.5
$ws.Range("M77").Value = -26276.75
$ws.Range("N77").Value = -17003.5
$ws.Range("H97").Value = 574.5
$ws.Range("I97").Value = 591.53845
$ws.Range("J97").Value = 542.8570999999999
$ws.Range("K97").Value = 591.53845
$ws.Range("L97").Value = 542.8570999999999
$ws.Range("M97").Value = -95.53845000000001
$ws.Range("N97").Value = -1534.8571
$ws.Range("H116").Value = 1370.3704
$ws.Range("I116").Value = 1207.1875
$ws.Range("J116").Value = 1607.7273
$ws.Range("K116").Value = 1207.1875
$ws.Range("L116").Value = 1607.7273
$ws.Range("M116").Value = 1086.8125
$ws.Range("N116").Value = -6195.7273
$ws.Range("H125").Value = 52185.69
$ws.Range("J125").Value = 52185.69
$ws.Range("L125").Value = 52185.69
$ws.Range("N125").Value = -62025.69
$ws.Range("H132").Value = 2042.8206
$ws.Range("I132").Value = 1883.027
$ws.Range("K132").Value = 5649.081
$ws.Range("M132").Value = -3119.081

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1370.3704
$ws.Range("I3").Value = 1207.1875
$ws.Range("J3").Value = 1607.7273
$ws.Range("K3").Value = 1207.1875
$ws.Range("L3").Value = 1607.7273
$ws.Range("M3").Value = -1093.1875
$ws.Range("N3").Value = -1835.7273
$ws.Range("H134").Value = 1426.2826
$ws.Range("I134").Value = 1426.2826
$ws.Range("K134").Value = 4278.8478
$ws.Range("M134").Value = -1743.8478

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3028.56
$ws.Range("I31").Value = 2538.3408
$ws.Range("K31").Value = 2538.3408
$ws.Range("M31").Value = -2243.3408
$ws.Range("H34").Value = 3028.56
$ws.Range("I34").Value = 2538.3408
$ws.Range("K34").Value = 2538.3408
$ws.Range("M34").Value = -2336.3408
$ws.Range("H58").Value = 1393.2122
$ws.Range("I58").Value = 631.4706
$ws.Range("K58").Value = 631.4706
$ws.Range("M58").Value = -428.4706
$ws.Range("H132").Value = 2002.5714
$ws.Range("I132").Value = 1441
$ws.Range("J132").Value = 3799.6
$ws.Range("K132").Value = 4323
$ws.Range("L132").Value = 11398.8
$ws.Range("M132").Value = -1793
$ws.Range("N132").Value = -16458.8
$ws.Range("H134").Value = 4501.8945
$ws.Range("I134").Value = 2533.625
$ws.Range("J134").Value = 14999.333
$ws.Range("K134").Value = 7600.875
$ws.Range("L134").Value = 44997.999
$ws.Range("M134").Value = -5065.875
$ws.Range("N134").Value = -50067.999
$ws.Range("H136").Value = 1393.2122
$ws.Range("I136").Value = 631.4706
$ws.Range("K136").Value = 1894.4118
$ws.Range("M136").Value = 655.5882000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 29.285715
$ws.Range("I2").Value = 29.285715
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 175.71429
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -62.71429000000001
$ws.Range("N2").ClearContents()
$ws.Range("H22").Value = 1854.2858
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1854.2858
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 5562.857400000001
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -5900.857400000001
$ws.Range("H27").Value = 1854.2858
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1854.2858
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 5562.857400000001
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -5766.857400000001
$ws.Range("H33").Value = 35.77778
$ws.Range("I33").Value = 30
$ws.Range("K33").Value = 180
$ws.Range("M33").Value = 103
$ws.Range("H44").Value = 1136
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 1136
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 3408
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -4204
$ws.Range("H92").Value = 5000500
$ws.Range("I92").Value = 10000000
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 30000000
$ws.Range("L92").Value = 3000
$ws.Range("M92").Value = -29998752
$ws.Range("N92").Value = -5496
$ws.Range("H94").Value = 3903.4285
$ws.Range("J94").Value = 4383.3335
$ws.Range("L94").Value = 13150.0005
$ws.Range("N94").Value = -14502.0005
$ws.Range("H132").Value = 674055.25
$ws.Range("I132").Value = 609.9167
$ws.Range("K132").Value = 5489.2503
$ws.Range("M132").Value = -2959.2503

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2335.4546
$ws.Range("I102").Value = 2598
$ws.Range("J102").Value = 2116.6667
$ws.Range("K102").Value = 2598
$ws.Range("L102").Value = 2116.6667
$ws.Range("M102").Value = -976
$ws.Range("N102").Value = -5360.6667
$ws.Range("H107").Value = 413.33334
$ws.Range("I107").Value = 270.1
$ws.Range("J107").Value = 699.8
$ws.Range("K107").Value = 270.1
$ws.Range("L107").Value = 699.8
$ws.Range("M107").Value = 1649.9
$ws.Range("N107").Value = -4539.8
$ws.Range("H113").Value = 10025.538
$ws.Range("I113").Value = 2030.1818
$ws.Range("K113").Value = 2030.1818
$ws.Range("M113").Value = 139.8181999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4177.923
$ws.Range("I40").Value = 2473.2856
$ws.Range("J40").Value = 6166.6665
$ws.Range("K40").Value = 2473.2856
$ws.Range("L40").Value = 6166.6665
$ws.Range("M40").Value = -2337.2856
$ws.Range("N40").Value = -6438.6665
$ws.Range("H68").Value = 3334666.8
$ws.Range("I68").Value = 3334666.8
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 3334666.8
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -3333917.8
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 3334666.8
$ws.Range("I71").Value = 3334666.8
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 16673334
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -16669590
$ws.Range("N71").Value = -16669590
$ws.Range("H132").Value = 5983.231
$ws.Range("I132").Value = 6911.125
$ws.Range("K132").Value = 20733.375
$ws.Range("M132").Value = -18203.375
$ws.Range("H136").Value = 1266.1
$ws.Range("I136").Value = 989.6667
$ws.Range("J136").Value = 1911.1111
$ws.Range("K136").Value = 2969.0001
$ws.Range("L136").Value = 5733.3333
$ws.Range("M136").Value = -419.0001000000002
$ws.Range("N136").Value = -10833.3333

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 66670
$ws.Range("J14").Value = 66670
$ws.Range("L14").Value = 66670
$ws.Range("N14").Value = -67006
$ws.Range("H39").Value = 5522
$ws.Range("I39").Value = 4377.3335
$ws.Range("J39").Value = 6666.6665
$ws.Range("K39").Value = 4377.3335
$ws.Range("L39").Value = 6666.6665
$ws.Range("M39").Value = -3964.3335
$ws.Range("N39").Value = -7492.6665
$ws.Range("H132").Value = 1787.6562
$ws.Range("I132").Value = 1092.3846
$ws.Range("K132").Value = 3277.1538
$ws.Range("M132").Value = -747.1538
$ws.Range("H136").Value = 4196.3687
$ws.Range("I136").Value = 4759.5713
$ws.Range("J136").Value = 2619.4
$ws.Range("K136").Value = 14278.7139
$ws.Range("L136").Value = 7858.200000000001
$ws.Range("M136").Value = -11728.7139
$ws.Range("N136").Value = -12958.2

Write-Output "Applied 223 value updates and 5 clears."